# "Menu & HUD af" - hours-registration workbook update
#
# 1) Clear every "Afwezigheid Omschrijving" (absence reason) note in column H
#    — the notes are gone, but the styling of the (now empty) cells stays.
# 2) A handful of day cells (columns C..G) get real attendance numbers filled
#    in / corrected. Every such cell is colour-coded: blue fill when the
#    person was present (value > 0), red fill when absent (value 0) — this
#    mirrors the colour convention already used throughout the rest of the
#    sheet.
# 3) Move the active selection to H4 (no more scrolled-down view).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

$blue = 12611584   # RGB(0,112,192)  -> present / value > 0
$red  = 255        # RGB(255,0,0)    -> absent  / value = 0

function Set-Attendance {
    param(
        [string]$addr,
        [double]$value
    )
    $rng = $ws.Range($addr)
    $rng.Value = $value
    if ($value -eq 0) {
        $rng.Interior.Color = $red
    } else {
        $rng.Interior.Color = $blue
    }
}

# --- 1. Clear the stale absence-reason comments in column H -----------------
$hCells = @(
    "H6","H7",
    "H11","H12","H13","H14","H15",
    "H19","H20","H21","H22","H23",
    "H42","H43","H44","H45","H46",
    "H50","H51","H53","H54",
    "H59","H60","H61","H62",
    "H66","H67","H68","H69","H70",
    "H74","H75","H76","H77"
)
foreach ($addr in $hCells) {
    $ws.Range($addr).Value = ""
}

# --- 2. Fix up / fill in attendance numbers ---------------------------------

# Week 4 (row 19-23): Fahrettin actually showed up on the Friday.
Set-Attendance "C23" 4

# Week 8 (row 42-46): Fahrettin actually showed up on the Monday.
Set-Attendance "C42" 4

# Week 10 (row 50-54): Fahrettin actually showed up on the Friday.
Set-Attendance "C53" 2

# Week 14 (rows 74-78): Friday row was still blank, now filled in.
Set-Attendance "C78" 0
Set-Attendance "D78" 4
Set-Attendance "E78" 4
Set-Attendance "F78" 4
Set-Attendance "G78" 4

# Week 15 (rows 82-86): Monday row was still blank, now filled in.
Set-Attendance "C82" 0
Set-Attendance "D82" 2
Set-Attendance "E82" 4
Set-Attendance "F82" 4
Set-Attendance "G82" 4

# Tuesday row already had numbers, just needed the colour coding applied.
Set-Attendance "C83" 4
Set-Attendance "D83" 0
Set-Attendance "E83" 4
Set-Attendance "F83" 4
Set-Attendance "G83" 4

# Wednesday row was still blank, now filled in.
Set-Attendance "C84" 2
Set-Attendance "D84" 0
Set-Attendance "E84" 2
Set-Attendance "F84" 2
Set-Attendance "G84" 2

# --- 3. Move the selection / scroll position --------------------------------
$ws.Range("H4").Select()
